$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: Wednesday (D), Thursday (E) slot at 8:40 - clear the MCT class
$ws.Range("E4").Value = "-"

# Row 6 (9:50 slot): shift the MCT-2A-MTRM / MEC-2A-MTRM classes
$ws.Range("B6").Value = "MCT-2A-MTRM"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "MEC-2A-MTRM"

# Row 7 (10:40 slot): shift the MCT-2A-MTRM / MEC-2A-MTRM classes
$ws.Range("B7").Value = "MCT-2A-MTRM"
$ws.Range("D7").Value = "-"
$ws.Range("F7").Value = "MEC-2A-MTRM"
